$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.474.92'
$ws.Range("E2").Value = '  +0.18%  '

# Row 3
$ws.Range("D3").Value = '3.500.66'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.77%  '

# Row 7
$ws.Range("D7").Value = '3.498.93'
$ws.Range("E7").Value = '  -0.49%  '

# Row 8
$ws.Range("E8").Value = '  +0.25%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.67%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.432'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.67%  '

# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000217'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.57%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.79%  '

# Row 15
$ws.Range("D15").Value = '4.093.07'
$ws.Range("E15").Value = '  -0.38%  '

# Row 16
$ws.Range("D16").Value = '3.503.33'
$ws.Range("E16").Value = '  -0.23%  '

# Row 17
$ws.Range("D17").Value = '67.428.96'
$ws.Range("E17").Value = '  +0.12%  '

# Row 18
$ws.Range("E18").Value = '  +0.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.15%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.93%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.631'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.96%  '

# Row 25
$ws.Range("D25").Value = '3.643.00'
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
$ws.Range("E26").Value = '  -0.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000126'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.47%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.70%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.30%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.65%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.170'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.59%  '

# Row 33
$ws.Range("E33").Value = '  +0.08%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.17'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.39%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.09%  '

# Row 37
$ws.Range("D37").Value = '3.496.11'
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("E39").Value = '  +0.00%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.45%  '

# Row 41
$ws.Range("E41").Value = '  +0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '175.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.29%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0889'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.32%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.883'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.29%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.20%  '

# Row 48
$ws.Range("E48").Value = '  +3.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.994'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.42%  '
